# "Generate Report for Archive" -- refresh the localization-status report:
# the handoff batch has moved from "Ready for handoff" into active
# translation, and the (now shorter) status text lets the date/status
# columns on each sheet be narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# 1. Status text: "Ready for handoff" -> "In Translation"
#    Overview!E2:F3 hold the per-locale status; zh-cn!C2:C3 / de-de!C2:C3
#    hold the same status on their own sheet.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# 2. Re-fit the now-narrower status/date columns to match the shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
